$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(1, 1).Value = 'Datos actualizados a 31 de Marzo de 2020 a las 15:50'

$ws.Cells.Item(20, 1).Value = 'Tenerife'
$ws.Cells.Item(20, 3).Value = 30
$ws.Cells.Item(20, 4).Value = 1056
$ws.Cells.Item(20, 5).Value = 36

$ws.Cells.Item(21, 1).Value = 'Gran Canaria'
$ws.Cells.Item(21, 4).Value = 342
$ws.Cells.Item(21, 5).Value = 13

$ws.Cells.Item(22, 1).Value = 'La Palma'
$ws.Cells.Item(22, 4).Value = 57
$ws.Cells.Item(22, 5).Value = 2

$ws.Cells.Item(23, 1).Value = 'Lanzarote'
$ws.Cells.Item(23, 4).Value = 45
$ws.Cells.Item(23, 5).Value = 3

$ws.Cells.Item(24, 1).Value = 'Fuerteventura'
$ws.Cells.Item(24, 4).Value = 32

$ws.Cells.Item(25, 1).Value = 'La Gomera'
$ws.Cells.Item(25, 4).Value = 7

$ws.Cells.Item(26, 1).Value = 'El Hierro'
$ws.Cells.Item(26, 2).Value = 1262
$ws.Cells.Item(26, 3).Value = 57
$ws.Cells.Item(26, 4).Value = 3
$ws.Cells.Item(26, 5).Value = 0

$ws.Cells.Item(27, 1).Value = 'Asturias'
$ws.Cells.Item(27, 2).Value = 1236
$ws.Cells.Item(27, 3).Value = 90
$ws.Cells.Item(27, 4).Value = 1091
$ws.Cells.Item(27, 5).Value = 55

$ws.Cells.Item(28, 1).Value = 'Salamanca'
$ws.Cells.Item(28, 2).Value = 1185
$ws.Cells.Item(28, 3).Value = 209
$ws.Cells.Item(28, 4).Value = 853
$ws.Cells.Item(28, 5).Value = 123

$ws.Cells.Item(29, 1).Value = 'Cantabria'
$ws.Cells.Item(29, 2).Value = 1171
$ws.Cells.Item(29, 3).Value = 35
$ws.Cells.Item(29, 4).Value = 1099
$ws.Cells.Item(29, 5).Value = 37

$ws.Cells.Item(30, 1).Value = 'Sevilla'
$ws.Cells.Item(30, 2).Value = 1119
$ws.Cells.Item(30, 3).Value = 17
$ws.Cells.Item(30, 4).Value = 1066
$ws.Cells.Item(30, 5).Value = 36

$ws.Cells.Item(31, 1).Value = 'Gipuzkoa/Guipuzcoa'
$ws.Cells.Item(31, 2).Value = 1113
$ws.Cells.Item(31, 3).Value = 1796
$ws.Cells.Item(31, 4).Value = 673
$ws.Cells.Item(31, 5).Value = 44

$ws.Cells.Item(32, 1).Value = 'Caceres'
$ws.Cells.Item(32, 2).Value = 1067
$ws.Cells.Item(32, 3).Value = 29
$ws.Cells.Item(32, 4).Value = 924
$ws.Cells.Item(32, 5).Value = 114

$ws.Cells.Item(33, 1).Value = 'Granada'
$ws.Cells.Item(33, 2).Value = 1061
$ws.Cells.Item(33, 3).Value = 15
$ws.Cells.Item(33, 4).Value = 979
$ws.Cells.Item(33, 5).Value = 67

$ws.Cells.Item(34, 1).Value = 'Valladolid'
$ws.Cells.Item(34, 2).Value = 988
$ws.Cells.Item(34, 3).Value = 187
$ws.Cells.Item(34, 4).Value = 724
$ws.Cells.Item(34, 5).Value = 77

$ws.Cells.Item(35, 1).Value = 'Leon'
$ws.Cells.Item(35, 2).Value = 964
$ws.Cells.Item(35, 3).Value = 176
$ws.Cells.Item(35, 4).Value = 688
$ws.Cells.Item(35, 5).Value = 100

$ws.Cells.Item(36, 1).Value = 'Murcia'
$ws.Cells.Item(36, 2).Value = 939
$ws.Cells.Item(36, 3).Value = 20
$ws.Cells.Item(36, 4).Value = 920
$ws.Cells.Item(36, 5).Value = 34

$ws.Cells.Item(37, 1).Value = 'Aragon'
$ws.Cells.Item(37, 2).Value = 907
$ws.Cells.Item(37, 3).Value = 29
$ws.Cells.Item(37, 4).Value = 838
$ws.Cells.Item(37, 5).Value = 40

$ws.Cells.Item(38, 1).Value = 'Burgos'
$ws.Cells.Item(38, 2).Value = 790
$ws.Cells.Item(38, 3).Value = 215
$ws.Cells.Item(38, 4).Value = 506
$ws.Cells.Item(38, 5).Value = 69

$ws.Cells.Item(39, 1).Value = 'Segovia'
$ws.Cells.Item(39, 2).Value = 720
$ws.Cells.Item(39, 3).Value = 189
$ws.Cells.Item(39, 4).Value = 454
$ws.Cells.Item(39, 5).Value = 77

$ws.Cells.Item(40, 1).Value = 'Cordoba'
$ws.Cells.Item(40, 3).Value = 4
$ws.Cells.Item(40, 4).Value = 642
$ws.Cells.Item(40, 5).Value = 15

$ws.Cells.Item(41, 1).Value = 'Jaen'
$ws.Cells.Item(41, 2).Value = 661
$ws.Cells.Item(41, 3).Value = 17
$ws.Cells.Item(41, 4).Value = 618
$ws.Cells.Item(41, 5).Value = 26

$ws.Cells.Item(42, 1).Value = 'Guadalajara'
$ws.Cells.Item(42, 2).Value = 643
$ws.Cells.Item(42, 3).Value = 296
$ws.Cells.Item(42, 4).Value = 532
$ws.Cells.Item(42, 5).Value = 97

$ws.Cells.Item(43, 1).Value = 'Castello/Castellon'
$ws.Cells.Item(43, 2).Value = 613
$ws.Cells.Item(43, 3).Value = 9
$ws.Cells.Item(43, 4).Value = 570
$ws.Cells.Item(43, 5).Value = 34

$ws.Cells.Item(44, 1).Value = 'Ourense'
$ws.Cells.Item(44, 2).Value = 570
$ws.Cells.Item(44, 3).Value = 187
$ws.Cells.Item(44, 4).Value = 520
$ws.Cells.Item(44, 5).Value = 12

$ws.Cells.Item(45, 1).Value = 'Badajoz'
$ws.Cells.Item(45, 2).Value = 561
$ws.Cells.Item(45, 3).Value = 62
$ws.Cells.Item(45, 4).Value = 480
$ws.Cells.Item(45, 5).Value = 19

$ws.Cells.Item(46, 1).Value = 'Soria'
$ws.Cells.Item(46, 2).Value = 555
$ws.Cells.Item(46, 3).Value = 87
$ws.Cells.Item(46, 4).Value = 424
$ws.Cells.Item(46, 5).Value = 44

$ws.Cells.Item(47, 1).Value = 'Cadiz'
$ws.Cells.Item(47, 2).Value = 539
$ws.Cells.Item(47, 3).Value = 16
$ws.Cells.Item(47, 4).Value = 509
$ws.Cells.Item(47, 5).Value = 14

$ws.Cells.Item(48, 1).Value = 'Avila'
$ws.Cells.Item(48, 2).Value = 467
$ws.Cells.Item(48, 3).Value = 111
$ws.Cells.Item(48, 4).Value = 302
$ws.Cells.Item(48, 5).Value = 54

$ws.Cells.Item(49, 1).Value = 'Lugo'
$ws.Cells.Item(49, 2).Value = 402
$ws.Cells.Item(49, 3).Value = 187
$ws.Cells.Item(49, 4).Value = 371
$ws.Cells.Item(49, 5).Value = 7

$ws.Cells.Item(50, 1).Value = 'Palencia'
$ws.Cells.Item(50, 2).Value = 325
$ws.Cells.Item(50, 3).Value = 42
$ws.Cells.Item(50, 4).Value = 265
$ws.Cells.Item(50, 5).Value = 18

$ws.Cells.Item(51, 1).Value = 'Cuenca'
$ws.Cells.Item(51, 2).Value = 293
$ws.Cells.Item(51, 3).Value = 296
$ws.Cells.Item(51, 4).Value = 210
$ws.Cells.Item(51, 5).Value = 64

$ws.Cells.Item(52, 1).Value = 'Almeria'
$ws.Cells.Item(52, 2).Value = 251
$ws.Cells.Item(52, 3).Value = 6
$ws.Cells.Item(52, 4).Value = 229
$ws.Cells.Item(52, 5).Value = 16

$ws.Cells.Item(53, 1).Value = 'Huesca'
$ws.Cells.Item(53, 2).Value = 244
$ws.Cells.Item(53, 3).Value = 23
$ws.Cells.Item(53, 4).Value = 207
$ws.Cells.Item(53, 5).Value = 14

$ws.Cells.Item(54, 1).Value = 'Teruel'
$ws.Cells.Item(54, 2).Value = 236
$ws.Cells.Item(54, 3).Value = 16
$ws.Cells.Item(54, 4).Value = 205
$ws.Cells.Item(54, 5).Value = 15

$ws.Cells.Item(55, 1).Value = 'Zamora'
$ws.Cells.Item(55, 2).Value = 217
$ws.Cells.Item(55, 3).Value = 43
$ws.Cells.Item(55, 4).Value = 151
$ws.Cells.Item(55, 5).Value = 23

$ws.Cells.Item(56, 1).Value = 'Mallorca'
$ws.Cells.Item(56, 2).Value = 210
$ws.Cells.Item(56, 3).Value = 18
$ws.Cells.Item(56, 4).Value = 194
$ws.Cells.Item(56, 5).Value = 12

$ws.Cells.Item(57, 1).Value = 'Huelva'
$ws.Cells.Item(57, 2).Value = 177
$ws.Cells.Item(57, 3).Value = 2
$ws.Cells.Item(57, 4).Value = 171
$ws.Cells.Item(57, 5).Value = 4

$ws.Cells.Item(58, 1).Value = 'Igualada, Vilanova del Cami, Santa Margarida de Montbui y Odena'
$ws.Cells.Item(58, 2).Value = 58
$ws.Cells.Item(58, 4).Value = 58
$ws.Cells.Item(58, 5).Value = 3

$ws.Cells.Item(59, 1).Value = 'Melilla'
$ws.Cells.Item(59, 2).Value = 54
$ws.Cells.Item(59, 4).Value = 53

$ws.Cells.Item(60, 1).Value = 'Ceuta'
$ws.Cells.Item(60, 2).Value = 29
$ws.Cells.Item(60, 3).Value = 0
$ws.Cells.Item(60, 4).Value = 28

$ws.Cells.Item(61, 1).Value = 'Ibiza'
$ws.Cells.Item(61, 2).Value = 21
$ws.Cells.Item(61, 4).Value = 20
$ws.Cells.Item(61, 5).Value = 1

$ws.Cells.Item(62, 1).Value = 'Menorca'
$ws.Cells.Item(62, 2).Value = 15
$ws.Cells.Item(62, 3).Value = 18
$ws.Cells.Item(62, 4).Value = 13

$ws.Cells.Item(63, 1).Value = 'Arroyo de la Luz'
$ws.Cells.Item(63, 2).Value = 7
$ws.Cells.Item(63, 3).Value = 0
$ws.Cells.Item(63, 4).Value = 7
$ws.Cells.Item(63, 5).Value = 0
